$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) stores free-form numeric-like text (e.g. "47.935.08", "1.00")
# that must survive verbatim -- trailing zeros, thousand-dot groupings, etc. --
# rather than being auto-coerced into a Double. Mark the specific Price cells
# we are about to rewrite as Text first so Excel keeps the literal string.
$priceRows = @(2,3,5,6,7,8,9,10,11,13,15,16,17,18,19,20,21,22,23,25,26,31,32,33,34,39,41,42,43,44,45,46,48,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "47.993.22"
$ws.Range("E2").Value = "  +6.30%  "
$ws.Range("D3").Value = "2.514.78"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "324.46"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").Value = "106.20"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").Value = "38.02"
$ws.Range("E10").Value = "  +7.13%  "
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "18.46"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "2.908.56"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "2.506.55"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "0.848"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "47.847.64"
$ws.Range("E18").Value = "  +6.13%  "
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  +4.42%  "
$ws.Range("D20").Value = "6.59"
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("D21").Value = "0.0₃0941"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "70.90"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("D23").Value = "251.67"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("E24").Value = "  +6.41%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").Value = "26.40"
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +5.06%  "
$ws.Range("E29").Value = "  +6.48%  "
$ws.Range("E30").Value = "  +7.06%  "
$ws.Range("D31").Value = "0.136"
$ws.Range("E31").Value = "  +9.73%  "
$ws.Range("D32").Value = "49.48"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "20.17"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "5.38"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "2.25"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "121.35"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("D43").Value = "21.32"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").Value = "0.0299"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("D45").Value = "1.972.52"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "3.01"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "1.83"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "9.24"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").Value = "5.40"
$ws.Range("E50").Value = "  +14.40%  "
$ws.Range("D51").Value = "79.38"
$ws.Range("E51").Value = "  +3.46%  "
